$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6041571497917175
$ws.Range("B1").Value = 1.198070406913757
$ws.Range("C1").Value = 5.037482261657715
$ws.Range("D1").Value = 1.443938851356506
$ws.Range("E1").Value = 0.8230966329574585
